$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Well column (C3:C12): previously held column-A well refs (A3..A12);
# now these are simple sequential column-B well refs for the plate (B3..B12).
# (C11 set before C10 to match the original authoring order of the shared
# string table.)
$ws.Range("C3").Value = "B3"
$ws.Range("C4").Value = "B4"
$ws.Range("C5").Value = "B5"
$ws.Range("C6").Value = "B6"
$ws.Range("C7").Value = "B7"
$ws.Range("C8").Value = "B8"
$ws.Range("C9").Value = "B9"
$ws.Range("C11").Value = "B11"
$ws.Range("C10").Value = "B10"
$ws.Range("C12").Value = "B12"

# --- Rows 10/11 had the pRC008.206 100x/1000x dilution label & concentration
# swapped; fix the ordering so row 10 is the 100x dilution and row 11 is the
# 1000x dilution.
$ws.Range("B10").Value = "pRC008.206_100x"
$ws.Range("D10").Value = 16.465
$ws.Range("B11").Value = "pRC008.206_1000x"
$ws.Range("D11").Value = 0.16465

# --- Water row (row 2): the Water source now feeds every well on the plate,
# so the well list / matching concentration list expand from the first two
# wells to the full A1:A12 set.
$ws.Range("C2").Value = "A1,A2,A3,A4,A5,A6,A7,A8,A9,A10,A11,A12"
$ws.Range("E2").Value = "65,65,65,65,65,65,65,65,65,65,65,65"

# --- Clear stray formatted-but-empty cells left over from earlier edits.
$ws.Range("F10:Y10").Clear()
$ws.Range("B13:E13").Clear()

# --- Leave the cursor where the author left it.
$ws.Range("F10").Select()
